# Apply the authored changes to sheet1 of the workbook:
#  - widen column A (to fit the newly added "luận giải" content)
#  - move the active selection down to A17 (below the data that was added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A from ~21.4 chars to ~36.9 chars.
$ws.Columns.Item(1).ColumnWidth = 36

# Move / leave the selection at A17.
$ws.Range("A17").Select()
